$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------- Step 1: add "evaluator_partial_correctness" header column to o_10 ----------
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("E1").Value = "evaluator_partial_correctness"

# ---------- Step 2: update o_10 row 2 data (prompt/solution/llm_response changed) ----------
$o10_prompt = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 10 nodes labelled A to J. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node J?
   A B C D E F G H I J
 A 0 1 0 0 0 0 0 0 0 0
 B 1 0 1 1 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0
 D 0 1 0 0 1 0 0 1 1 0
 E 0 0 0 1 0 1 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0
 G 0 0 0 0 0 1 0 0 0 0
 H 0 0 0 1 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 0 0 1
 J 0 0 0 0 0 0 0 0 1 0
    
"@
$o10_solution = @"
A -> B -> D -> I -> J
"@
$o10_llm_response = @"
The shortest path from node A to node J is A -> B -> D -> H -> J.
"@
$o10_eval_partial = @"
Output: 3/5
"@

$ws1.Range("A2").Value = $o10_prompt
$ws1.Range("B2").Value = $o10_solution
$ws1.Range("C2").Value = $o10_llm_response
$ws1.Range("E2").Value = $o10_eval_partial
$ws1.Rows.Item(2).EntireRow.AutoFit()

# ---------- Step 3: add new sheets o_20 and o_20_jumbled (after o_10) ----------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "o_20"
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "o_20_jumbled"

# ---------- Step 4: copy header row (with formatting) to new sheets ----------
$ws1.Range("A1:E1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)
$ws1.Range("A1:E1").Copy()
$ws3.Range("A1:E1").PasteSpecial(-4122)

$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"

$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"

# ---------- Step 5: fill o_20 row 2 data ----------
$o20_prompt = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node T?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 1 0 0 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 1 0 0 1 1 1 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 1 0 0 0 1 1 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 1 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 1 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
"@
$o20_solution = @"
A -> B -> E -> G -> J -> L -> M -> N -> P -> Q -> S -> T
"@
$o20_llm_response = @"
To find the shortest path from node A to node T, we can use the Breadth-First Search (BFS) algorithm. 
1. Initialize a queue and an empty list that will store the visited nodes.
2. Enqueue node A into the queue.
3. Mark node A as visited by adding it to the visited list.
4. While the queue is not empty, repeat steps 5-8.
5. Dequeue a node from the queue.
6. Check if the dequeued node is node T. If it is, we have found the shortest path from node A to node T. 
7. If the dequeued node is not node T, find all its adjacent nodes that have not been visited yet.
8. Enqueue each unvisited adjacent node into the queue and mark them as visited by adding them to the visited list.
9. If the queue becomes empty, there is no path from node A to node T.
By following this algorithm, we can find the shortest path from node A to node T.
"@
$o20_eval_response = "Wrong"
$o20_eval_partial = @"
Output: 0/0
"@

$ws2.Range("A2").Value = $o20_prompt
$ws2.Range("B2").Value = $o20_solution
$ws2.Range("C2").Value = $o20_llm_response
$ws2.Range("D2").Value = $o20_eval_response
$ws2.Range("E2").Value = $o20_eval_partial
$ws2.Rows.Item(2).EntireRow.AutoFit()

# ---------- Step 6: fill o_20_jumbled row 2 data ----------
$o20j_prompt = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node T?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 1 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 1 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 1
 K 0 0 0 0 0 0 0 0 0 1 0 1 1 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 1 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
    
"@
$o20j_solution = @"
A -> I -> J -> T
"@
$o20j_llm_response = @"
To find the shortest path from node A to node T, we can use a breadth-first search algorithm. 
Starting from node A, we can explore its neighboring nodes (B, C, D, F, and I). From these nodes, we can explore their neighboring nodes and so on, until we reach node T.
We can maintain a queue to keep track of the nodes to explore next. Initially, we add node A to the queue. 
Then, while the queue is not empty, we take out the first node from the queue and explore its neighbors. 
In this case, the neighbors of node A are B, C, D, F, and I. We can add these neighbors to the queue and mark them as visited. 
Next, we take out the next node from the queue, which is node B, and explore its neighbors. Node B has no neighbors in this case, so we move on to the next node in the queue. 
We continue this process until we reach node T or until the queue becomes empty. If we reach node T, we have found the shortest path from node A to node T. 
In the given adjacency matrix, we can observe that nodes A and T are not directly connected, but there is a path that connects them indirectly through other nodes.
Here is the step-by-step process:
1. Add node A to the queue.
2. Remove node A from the queue and mark it as visited.
3. Explore the neighbors of node A (B, C, D, F, and I).
4. Add the unvisited neighbors (B, C, D, F, and I) to the queue and mark them as visited.
5. Remove the next node from the queue.
6. Continue this process until node T is reached or the queue becomes empty.
7. If node T is reached, we have found the shortest path from node A to node T.
8. If the queue becomes empty before reaching node T, there is no path from node A to node T.
Note: Since the graph is unweighted, we only need to count the number of edges to get the shortest path. We do not consider the weights of the edges.
"@
$o20j_eval_response = "Wrong"
$o20j_eval_partial = @"
Output: 1/3
"@

$ws3.Range("A2").Value = $o20j_prompt
$ws3.Range("B2").Value = $o20j_solution
$ws3.Range("C2").Value = $o20j_llm_response
$ws3.Range("D2").Value = $o20j_eval_response
$ws3.Range("E2").Value = $o20j_eval_partial
$ws3.Rows.Item(2).EntireRow.AutoFit()

# ---------- Step 7: restore active sheet/tab selection to o_10 ----------
$ws1.Activate()
$ws1.Range("A1").Select()

Write-Host "Edit complete"
